$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.282155871391296
$ws.Range("B1").Value = 2.256068468093872
$ws.Range("C1").Value = 4.74271821975708
$ws.Range("D1").Value = 3.05150318145752
$ws.Range("E1").Value = 1.345105290412903
